$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bold/bordered/centered style used on A2 down through A3:A21
$ws.Range("A2").Copy()
$ws.Range("A3:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Flattened data for rows 2..21, 6 values per row: A,B,C,D,E,F
# (A=period, B=participant, C=dwg, D=participantAvgHullArea-n, E=participantAvgHullArea, F=participantTime)
$data = @(
    0, 3000, 1, 10, 6.60869307700319, 1070.6317,
    1, 3000, 2, 10, 7.299787708709571, 1415.9084,
    2, 3000, 3, 10, 7.85893764337613, 642.5690000000001,
    3, 3000, 4, 10, 4.593158806225821, 877.062,
    4, 3000, 5, 10, 5.17807708575072, 728.8030999999999,
    5, 3000, 6, 10, 4.725605545853591, 873.9599999999998,
    6, 3000, 7, 10, 5.333819297697558, 1330.8975,
    7, 3000, 8, 10, 5.108525760744161, 600.6339999999999,
    8, 3000, 9, 10, 7.377457297111516, 755.0204,
    9, 3000, 10, 10, 6.738185545039576, 1522.4909,
    10, 3000, 11, 10, 6.018201625009492, 1863.4546,
    11, 3000, 12, 10, 7.005961915086409, 255.7185,
    12, 3000, 13, 10, 5.040659110120907, 1565.5031,
    13, 3000, 14, 10, 3.597975746048344, 1179.2965,
    14, 3000, 15, 10, 4.420036612820993, 803.2319000000001,
    15, 3000, 16, 10, 4.464905601109714, 800.5603,
    16, 3000, 17, 10, 5.116316240997665, 1236.5861,
    17, 3000, 18, 10, 2.906358479492045, 1097.7146,
    18, 3000, 19, 10, 4.777966849038944, 1559.909399999999,
    19, 3000, 20, 10, 4.387744588261096, 779.5659999999999
)

$cols = 6
$startRow = 2
$rowCount = 20
for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $data[$i * $cols + $j]
    }
}
